# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on the
# zh-cn and de-de report sheets (rows 4 and 5) to reflect a
# newly generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-20 16:20:50"
$wsZhCn.Range("E5").Value = "2016-03-20 16:20:50"
$wsZhCn.Range("H4").Value = "2016-03-20 16:21:11"
$wsZhCn.Range("H5").Value = "2016-03-20 16:21:11"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-20 16:20:53"
$wsDeDe.Range("E5").Value = "2016-03-20 16:20:53"
$wsDeDe.Range("H4").Value = "2016-03-20 16:21:17"
$wsDeDe.Range("H5").Value = "2016-03-20 16:21:17"
